$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header style (same as H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I (I0) and J (IF), rows 2-13
$values = @{
    2  = @(1, 5)
    3  = @(7, 9)
    4  = @(4, 6)
    5  = @(1, 4)
    6  = @(1, 3)
    7  = @(1, 5)
    8  = @(1, 4)
    9  = @(1, 5)
    10 = @(1, 5)
    11 = @(7, 9)
    12 = @(7, 8)
    13 = @(5, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
